# Apply odds updates to "Jogos da Semana" sheet as described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.14
$ws.Range("K2").Value = 5.5

# Row 4
$ws.Range("G4").Value = 2.22
$ws.Range("H4").Value = 2.95
$ws.Range("I4").Value = 3.45
$ws.Range("K4").Value = 5.6
$ws.Range("N4").Value = 2.37
$ws.Range("O4").Value = 1.52
$ws.Range("P4").Value = 1.52
$ws.Range("Q4").Value = 2.37
$ws.Range("U4").Value = 9.75
$ws.Range("W4").Value = 22
$ws.Range("Z4").Value = 5.6
$ws.Range("AA4").Value = 5.8
$ws.Range("AB4").Value = 17
$ws.Range("AD4").Value = 7.7
$ws.Range("AE4").Value = 16.5
$ws.Range("AF4").Value = 12.5
$ws.Range("AH4").Value = 40
$ws.Range("AI4").Value = 55

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 3.9
$ws.Range("L7").Value = 1.29
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 1.85
$ws.Range("O7").Value = 1.75
$ws.Range("P7").Value = 1.4
$ws.Range("Q7").Value = 2.52
$ws.Range("R7").Value = 1.75
$ws.Range("S7").Value = 1.85
$ws.Range("T7").Value = 6.9
$ws.Range("V7").Value = 8.25
$ws.Range("W7").Value = 15.5
$ws.Range("X7").Value = 15
$ws.Range("Y7").Value = 27
$ws.Range("Z7").Value = 9.75
$ws.Range("AA7").Value = 6.6
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 70
$ws.Range("AD7").Value = 11.25
$ws.Range("AE7").Value = 22
$ws.Range("AH7").Value = 35
$ws.Range("AI7").Value = 40
$ws.Range("AJ7").Value = 600

# Row 17
$ws.Range("G17").Value = 3.1
$ws.Range("I17").Value = 2.3
$ws.Range("K17").Value = 7.5
$ws.Range("T17").Value = 8
$ws.Range("U17").Value = 15
$ws.Range("V17").Value = 12
$ws.Range("W17").Value = 34
$ws.Range("X17").Value = 29
$ws.Range("AD17").Value = 6.5
$ws.Range("AE17").Value = 10
$ws.Range("AG17").Value = 21
$ws.Range("AH17").Value = 21
